$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes one cell whose inline-string text content changed.
# "Numeric" entries hold a value that Excel would otherwise auto-convert
# to a real number (e.g. "9.99"), so those cells get an explicit text
# number format ("@") applied first, to keep them stored as text -
# matching the original workbook's t="inlineStr" / text cells.
$changes = @(
    @{ Cell = "D2"; Value = "69.126.53"; Numeric = $false },
    @{ Cell = "E2"; Value = "  +1.36%  "; Numeric = $false },
    @{ Cell = "D3"; Value = "3.778.37"; Numeric = $false },
    @{ Cell = "E3"; Value = "  -0.56%  "; Numeric = $false },
    @{ Cell = "E4"; Value = "  +0.34%  "; Numeric = $false },
    @{ Cell = "D5"; Value = "627.94"; Numeric = $true },
    @{ Cell = "E5"; Value = "  +4.48%  "; Numeric = $false },
    @{ Cell = "D6"; Value = "163.69"; Numeric = $true },
    @{ Cell = "E6"; Value = "  -0.83%  "; Numeric = $false },
    @{ Cell = "D7"; Value = "3.775.68"; Numeric = $false },
    @{ Cell = "E7"; Value = "  -0.50%  "; Numeric = $false },
    @{ Cell = "E8"; Value = "  +0.01%  "; Numeric = $false },
    @{ Cell = "D9"; Value = "0.518"; Numeric = $true },
    @{ Cell = "E9"; Value = "  +0.01%  "; Numeric = $false },
    @{ Cell = "E10"; Value = "  +0.55%  "; Numeric = $false },
    @{ Cell = "D11"; Value = "0.451"; Numeric = $true },
    @{ Cell = "E11"; Value = "  -0.22%  "; Numeric = $false },
    @{ Cell = "D12"; Value = "6.61"; Numeric = $true },
    @{ Cell = "E12"; Value = "  +2.28%  "; Numeric = $false },
    @{ Cell = "D13"; Value = "0.0000247"; Numeric = $true },
    @{ Cell = "E13"; Value = "  -1.12%  "; Numeric = $false },
    @{ Cell = "D14"; Value = "35.42"; Numeric = $true },
    @{ Cell = "E14"; Value = "  -1.12%  "; Numeric = $false },
    @{ Cell = "D15"; Value = "4.427.17"; Numeric = $false },
    @{ Cell = "E15"; Value = "  -0.26%  "; Numeric = $false },
    @{ Cell = "D16"; Value = "3.889.10"; Numeric = $false },
    @{ Cell = "E16"; Value = "  +2.33%  "; Numeric = $false },
    @{ Cell = "D17"; Value = "69.173.92"; Numeric = $false },
    @{ Cell = "E17"; Value = "  +1.36%  "; Numeric = $false },
    @{ Cell = "D18"; Value = "17.93"; Numeric = $true },
    @{ Cell = "E18"; Value = "  -2.79%  "; Numeric = $false },
    @{ Cell = "E19"; Value = "  -0.97%  "; Numeric = $false },
    @{ Cell = "D20"; Value = "7.06"; Numeric = $true },
    @{ Cell = "E20"; Value = "  -0.51%  "; Numeric = $false },
    @{ Cell = "D21"; Value = "467.28"; Numeric = $true },
    @{ Cell = "E21"; Value = "  +1.16%  "; Numeric = $false },
    @{ Cell = "D22"; Value = "9.59"; Numeric = $true },
    @{ Cell = "E22"; Value = "  -1.21%  "; Numeric = $false },
    @{ Cell = "D23"; Value = "0.702"; Numeric = $true },
    @{ Cell = "E24"; Value = "  -0.71%  "; Numeric = $false },
    @{ Cell = "D25"; Value = "83.16"; Numeric = $true },
    @{ Cell = "E25"; Value = "  +0.16%  "; Numeric = $false },
    @{ Cell = "D26"; Value = "12.01"; Numeric = $true },
    @{ Cell = "E26"; Value = "  -0.46%  "; Numeric = $false },
    @{ Cell = "E27"; Value = "  +1.44%  "; Numeric = $false },
    @{ Cell = "E28"; Value = "  +0.05%  "; Numeric = $false },
    @{ Cell = "D29"; Value = "9.99"; Numeric = $true },
    @{ Cell = "E29"; Value = "  -0.04%  "; Numeric = $false },
    @{ Cell = "D30"; Value = "3.942.34"; Numeric = $false },
    @{ Cell = "E30"; Value = "  -0.18%  "; Numeric = $false },
    @{ Cell = "D31"; Value = "2.66"; Numeric = $true },
    @{ Cell = "E31"; Value = "  +0.39%  "; Numeric = $false },
    @{ Cell = "E32"; Value = "  -1.00%  "; Numeric = $false },
    @{ Cell = "D33"; Value = "7.24"; Numeric = $true },
    @{ Cell = "E33"; Value = "  -1.39%  "; Numeric = $false },
    @{ Cell = "D34"; Value = "28.85"; Numeric = $true },
    @{ Cell = "E34"; Value = "  -1.91%  "; Numeric = $false },
    @{ Cell = "E35"; Value = "  +0.02%  "; Numeric = $false },
    @{ Cell = "D36"; Value = "3.726.53"; Numeric = $false },
    @{ Cell = "E36"; Value = "  -0.60%  "; Numeric = $false },
    @{ Cell = "D37"; Value = "8.96"; Numeric = $true },
    @{ Cell = "E37"; Value = "  -0.75%  "; Numeric = $false },
    @{ Cell = "E38"; Value = "  +2.85%  "; Numeric = $false },
    @{ Cell = "D39"; Value = "0.150"; Numeric = $true },
    @{ Cell = "E39"; Value = "  +8.11%  "; Numeric = $false },
    @{ Cell = "E40"; Value = "  -0.22%  "; Numeric = $false },
    @{ Cell = "D41"; Value = "5.84"; Numeric = $true },
    @{ Cell = "E41"; Value = "  +0.33%  "; Numeric = $false },
    @{ Cell = "E42"; Value = "  -1.98%  "; Numeric = $false },
    @{ Cell = "E43"; Value = "  +0.01%  "; Numeric = $false },
    @{ Cell = "D45"; Value = "154.94"; Numeric = $true },
    @{ Cell = "E45"; Value = "  +2.06%  "; Numeric = $false },
    @{ Cell = "D46"; Value = "0.297"; Numeric = $true },
    @{ Cell = "E46"; Value = "  -1.10%  "; Numeric = $false },
    @{ Cell = "B47"; Value = "OKB"; Numeric = $false },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; Numeric = $false },
    @{ Cell = "D47"; Value = "46.84"; Numeric = $true },
    @{ Cell = "E47"; Value = "  -1.30%  "; Numeric = $false },
    @{ Cell = "B48"; Value = "Stacks"; Numeric = $false },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; Numeric = $false },
    @{ Cell = "D48"; Value = "1.92"; Numeric = $true },
    @{ Cell = "E48"; Value = "  +2.29%  "; Numeric = $false },
    @{ Cell = "D49"; Value = "42.41"; Numeric = $true },
    @{ Cell = "E49"; Value = "  -1.66%  "; Numeric = $false },
    @{ Cell = "B50"; Value = "Cosmos"; Numeric = $false },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; Numeric = $false },
    @{ Cell = "D50"; Value = "8.40"; Numeric = $true },
    @{ Cell = "E50"; Value = "  +0.27%  "; Numeric = $false },
    @{ Cell = "B51"; Value = "ONDO"; Numeric = $false },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"; Numeric = $false },
    @{ Cell = "D51"; Value = "1.38"; Numeric = $true },
    @{ Cell = "E51"; Value = "  +2.03%  "; Numeric = $false }
)

foreach ($change in $changes) {
    $rng = $ws.Range($change.Cell)
    if ($change.Numeric) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $change.Value
}
